# Opportunities & Engagement - 1st Merge - 12 Mar 2025

$wb = $excel.ActiveWorkbook

# --- AddOpportunity sheet: update a few data cells ---
$ws = $wb.Worksheets.Item("AddOpportunity")

$ws.Range("C2").Value = "Buyside"
$ws.Range("D2").Value = "HC - Healthcare"
$ws.Range("E2").Value = "Dental"

$ws.Range("E3").Value = "CSDN-0000001546"
$ws.Range("E3").VerticalAlignment = -4108
$ws.Range("E3").WrapText = $true

# Copy E3's freshly-built style onto E4 (reuses the same style record instead
# of forking a second, near-duplicate one), then set E4's value.
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Value = "CSDN-0000001546"

# --- Users sheet: swap referral staff name ---
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "James Craven"
$usersSheet.Range("B4").Select() | Out-Null

# --- Contact sheet: move selection ---
$contactSheet = $wb.Worksheets.Item("Contact")
$contactSheet.Range("K8").Select() | Out-Null

# --- Restore AddOpportunity as the active sheet/selection ---
$ws.Activate() | Out-Null
$ws.Range("E2").Select() | Out-Null
